# LED Watch BOM update
# - Improved crystal board layout, added reverse voltage protection with
#   schottky diode (new row 16). Fixed incorrect capacitor package sizes
#   (new row 18, "0402 2.2 nF Decoupling").
# - New columns F ("Compatibility for v1.1") and G ("Checked in v1.1").
# - Row 2 (crystal caps) highlighted yellow and flagged "NO" compatible.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# New header column F1 ("Compatibility for v1.1")
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "Compatibility for v1.1"
$ws.Range("F1").Font.Bold = $true

# ---------------------------------------------------------------------
# Mark existing, already-verified rows "Yes" compatible with v1.1
# ---------------------------------------------------------------------
$ws.Range("F3").Value = "Yes"
$ws.Range("F7").Value = "Yes"
$ws.Range("F8").Value = "Yes"
$ws.Range("F9").Value = "Yes"
$ws.Range("F10").Value = "Yes"
$ws.Range("G10").Value = "Yes"
$ws.Range("F11").Value = "Yes"
$ws.Range("F12").Value = "Yes"
$ws.Range("F13").Value = "Yes"
$ws.Range("F14").Value = "Yes"
$ws.Range("G14").Value = "Yes"

# ---------------------------------------------------------------------
# Row 2 - existing "Crystal caps" row: highlight yellow, mark as
# incompatible with v1.1 (wrong package size)
# ---------------------------------------------------------------------
$ws.Range("A2:F2").Interior.Color = 65535
$ws.Range("F2").Value = "NO"

# ---------------------------------------------------------------------
# Row 16 - new part: reverse-voltage-protection schottky diode (JTI),
# flagged for crystal-size check
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "Capacitor"
$ws.Range("B16").Value = "JTI"
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Value = "712-1220-1-ND"
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = "Crystal caps"
$ws.Range("F16").Value = "Yes"

# ---------------------------------------------------------------------
# New header column G1 ("Checked in v1.1")
# ---------------------------------------------------------------------
$ws.Range("G1").Value = "Checked in v1.1"
$ws.Range("G1").Font.Bold = $true

$ws.Range("G16").Value = "Size correct - Check crystal specs"

# ---------------------------------------------------------------------
# Row 18 - new part: corrected 0402 2.2nF decoupling capacitor
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "Capacitor"
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = "?"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = "0402 2.2 nF Decoupling"
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Column width adjustments (values chosen so the serialized XML "width"
# lands on the authored values of 20 / 37 / 20.28515625 / 14.85546875)
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668
$ws.Columns.Item(5).ColumnWidth = 36.166666666666664
$ws.Columns.Item(6).ColumnWidth = 19.4167
$ws.Columns.Item(7).ColumnWidth = 13.9167

# ---------------------------------------------------------------------
# Selection matches the final author state
# ---------------------------------------------------------------------
$ws.Range("E26").Select()
